$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efna1"
$ws.Range("C2").Value = "Epha4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2.0
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 21.78783233333333
$ws.Range("H2").Value = 65.363497
$ws.Range("I2").Value = 0.9024488799587679
$ws.Range("J2").Value = 0.9024488799587679
$ws.Range("K2").Value = 2.0
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 6.708176333333333
$ws.Range("N2").Value = 20.124529
$ws.Range("O2").Value = 0.4356329228871633
$ws.Range("P2").Value = 0.4356329228871633
$ws.Range("Q2").Value = 146.1566212131014
$ws.Range("R2").Value = 1315.409590917913
$ws.Range("S2").Value = 0.3931364433326848
$ws.Range("T2").Value = 0.3931364433326848

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efna1"
$ws.Range("C3").Value = "Epha4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2.0
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 21.78783233333333
$ws.Range("H3").Value = 65.363497
$ws.Range("I3").Value = 0.9024488799587679
$ws.Range("J3").Value = 0.9024488799587679
$ws.Range("K3").Value = 3.0
$ws.Range("L3").Value = 1.0
$ws.Range("M3").Value = 6.789877333333333
$ws.Range("N3").Value = 20.369632
$ws.Range("O3").Value = 0.4409386339573907
$ws.Range("P3").Value = 0.4409386339573907
$ws.Range("Q3").Value = 147.9367089025671
$ws.Range("R3").Value = 1331.430380123104
$ws.Range("S3").Value = 0.3979245763453964
$ws.Range("T3").Value = 0.3979245763453964

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Efna1"
$ws.Range("C4").Value = "Epha4"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2.0
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 21.78783233333333
$ws.Range("H4").Value = 65.363497
$ws.Range("I4").Value = 0.9024488799587679
$ws.Range("J4").Value = 0.9024488799587679
$ws.Range("K4").Value = 3.0
$ws.Range("L4").Value = 1.0
$ws.Range("M4").Value = 1.900636333333334
$ws.Range("N4").Value = 5.701909000000001
$ws.Range("O4").Value = 0.1234284431554459
$ws.Range("P4").Value = 0.1234284431554459
$ws.Range("Q4").Value = 41.41074575730811
$ws.Range("R4").Value = 372.696711815773
$ws.Range("S4").Value = 0.1113878602806866
$ws.Range("T4").Value = 0.1113878602806866

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Efna1"
$ws.Range("C5").Value = "Epha4"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3.0
$ws.Range("F5").Value = 1.0
$ws.Range("G5").Value = 1.359006333333333
$ws.Range("H5").Value = 4.077019
$ws.Range("I5").Value = 0.05628984676448105
$ws.Range("J5").Value = 0.05628984676448104
$ws.Range("K5").Value = 2.0
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 6.708176333333333
$ws.Range("N5").Value = 20.124529
$ws.Range("O5").Value = 0.4356329228871633
$ws.Range("P5").Value = 0.4356329228871633
$ws.Range("Q5").Value = 9.116454122116778
$ws.Range("R5").Value = 82.04808709905099
$ws.Range("S5").Value = 0.02452171047488141
$ws.Range("T5").Value = 0.02452171047488141

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Efna1"
$ws.Range("C6").Value = "Epha4"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3.0
$ws.Range("F6").Value = 1.0
$ws.Range("G6").Value = 1.359006333333333
$ws.Range("H6").Value = 4.077019
$ws.Range("I6").Value = 0.05628984676448105
$ws.Range("J6").Value = 0.05628984676448104
$ws.Range("K6").Value = 3.0
$ws.Range("L6").Value = 1.0
$ws.Range("M6").Value = 6.789877333333333
$ws.Range("N6").Value = 20.369632
$ws.Range("O6").Value = 0.4409386339573907
$ws.Range("P6").Value = 0.4409386339573907
$ws.Range("Q6").Value = 9.227486298556444
$ws.Range("R6").Value = 83.047376687008
$ws.Range("S6").Value = 0.02482036813800112
$ws.Range("T6").Value = 0.02482036813800112

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Efna1"
$ws.Range("C7").Value = "Epha4"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3.0
$ws.Range("F7").Value = 1.0
$ws.Range("G7").Value = 1.359006333333333
$ws.Range("H7").Value = 4.077019
$ws.Range("I7").Value = 0.05628984676448105
$ws.Range("J7").Value = 0.05628984676448104
$ws.Range("K7").Value = 3.0
$ws.Range("L7").Value = 1.0
$ws.Range("M7").Value = 1.900636333333334
$ws.Range("N7").Value = 5.701909000000001
$ws.Range("O7").Value = 0.1234284431554459
$ws.Range("P7").Value = 0.1234284431554459
$ws.Range("Q7").Value = 2.582976814363445
$ws.Range("R7").Value = 23.246791329271
$ws.Range("S7").Value = 0.00694776815159851
$ws.Range("T7").Value = 0.006947768151598509

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Efna1"
$ws.Range("C8").Value = "Epha4"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3.0
$ws.Range("F8").Value = 1.0
$ws.Range("G8").Value = 0.9961713333333334
$ws.Range("H8").Value = 2.988514
$ws.Range("I8").Value = 0.04126127327675106
$ws.Range("J8").Value = 0.04126127327675105
$ws.Range("K8").Value = 2.0
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 6.708176333333333
$ws.Range("N8").Value = 20.124529
$ws.Range("O8").Value = 0.4356329228871633
$ws.Range("P8").Value = 0.4356329228871633
$ws.Range("Q8").Value = 6.682492962211779
$ws.Range("R8").Value = 60.142436659906
$ws.Range("S8").Value = 0.01797476907959707
$ws.Range("T8").Value = 0.01797476907959706

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Efna1"
$ws.Range("C9").Value = "Epha4"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3.0
$ws.Range("F9").Value = 1.0
$ws.Range("G9").Value = 0.9961713333333334
$ws.Range("H9").Value = 2.988514
$ws.Range("I9").Value = 0.04126127327675106
$ws.Range("J9").Value = 0.04126127327675105
$ws.Range("K9").Value = 3.0
$ws.Range("L9").Value = 1.0
$ws.Range("M9").Value = 6.789877333333333
$ws.Range("N9").Value = 20.369632
$ws.Range("O9").Value = 0.4409386339573907
$ws.Range("P9").Value = 0.4409386339573907
$ws.Range("Q9").Value = 6.763881156316445
$ws.Range("R9").Value = 60.874930406848
$ws.Range("S9").Value = 0.0181936894739932
$ws.Range("T9").Value = 0.0181936894739932

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Efna1"
$ws.Range("C10").Value = "Epha4"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3.0
$ws.Range("F10").Value = 1.0
$ws.Range("G10").Value = 0.9961713333333334
$ws.Range("H10").Value = 2.988514
$ws.Range("I10").Value = 0.04126127327675106
$ws.Range("J10").Value = 0.04126127327675105
$ws.Range("K10").Value = 3.0
$ws.Range("L10").Value = 1.0
$ws.Range("M10").Value = 1.900636333333334
$ws.Range("N10").Value = 5.701909000000001
$ws.Range("O10").Value = 0.1234284431554459
$ws.Range("P10").Value = 0.1234284431554459
$ws.Range("Q10").Value = 1.893359430358445
$ws.Range("R10").Value = 17.04023487322601
$ws.Range("S10").Value = 0.005092814723160788
$ws.Range("T10").Value = 0.005092814723160788
